$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New day of data appended to the log (row 48): date, total count,
# session-timeout errors, errors requiring analysis.
$ws.Range("A48").Value = 46007
$ws.Range("B48").Value = 615
$ws.Range("C48").Value = 22
$ws.Range("D48").Value = 593

# Scroll the view down to keep the newly entered row in sight and move the
# selection onto it, matching where the author was working.
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A48:D48").Select()
